$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 64479.938
$ws.Range("J28").Value = 3657.7144
$ws.Range("L28").Value = 3657.7144
$ws.Range("N28").Value = -4627.7144

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 17500
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4007.0833
$ws.Range("J137").Value = 4630.5835
$ws.Range("L137").Value = 13891.7505
$ws.Range("N137").Value = -18991.7505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3328.4517
$ws.Range("I45").Value = 1706.5555
$ws.Range("J45").Value = 5574.154
$ws.Range("K45").Value = 1706.5555
$ws.Range("L45").Value = 5574.154
$ws.Range("M45").Value = -1329.5555
$ws.Range("N45").Value = -6328.154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3622.4849
$ws.Range("I61").Value = 2506
$ws.Range("J61").Value = 6190.4
$ws.Range("K61").Value = 2506
$ws.Range("L61").Value = 6190.4
$ws.Range("M61").Value = -2294
$ws.Range("N61").Value = -6614.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2743.75
$ws.Range("I74").Value = 2566.261
$ws.Range("K74").Value = 2566.261
$ws.Range("M74").Value = -1692.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2743.75
$ws.Range("I77").Value = 2566.261
$ws.Range("K77").Value = 12831.305
$ws.Range("M77").Value = -8463.305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 12500
$ws.Range("J94").Value = 12500
$ws.Range("L94").Value = 12500
$ws.Range("N94").Value = -14302

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3622.4849
$ws.Range("I136").Value = 2506
$ws.Range("J136").Value = 6190.4
$ws.Range("K136").Value = 7518
$ws.Range("L136").Value = 18571.2
$ws.Range("M136").Value = -4968
$ws.Range("N136").Value = -23671.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 500
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 588058.9399999999
$ws.Range("I86").Value = 1701914.5
$ws.Range("J86").Value = 1819.1578
$ws.Range("K86").Value = 1701914.5
$ws.Range("L86").Value = 1819.1578
$ws.Range("M86").Value = -1700791.5
$ws.Range("N86").Value = -4065.1578

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 588058.9399999999
$ws.Range("I89").Value = 1701914.5
$ws.Range("J89").Value = 1819.1578
$ws.Range("K89").Value = 8509572.5
$ws.Range("L89").Value = 9095.789000000001
$ws.Range("M89").Value = -8503956.5
$ws.Range("N89").Value = -20327.789

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2166.6667
$ws.Range("I105").Value = 1500
$ws.Range("K105").Value = 1500
$ws.Range("M105").Value = 247

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 1671.4286
$ws.Range("I32").Value = 360
$ws.Range("J32").Value = 4950
$ws.Range("K32").Value = 360
$ws.Range("L32").Value = 4950
$ws.Range("M32").Value = -44
$ws.Range("N32").Value = -5582

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 882
$ws.Range("J35").Value = 1000
$ws.Range("L35").Value = 1000
$ws.Range("N35").Value = -1588

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 11013
$ws.Range("I17").Value = 11013
$ws.Range("K17").Value = 33039
$ws.Range("M17").Value = -32870

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 399.5
$ws.Range("J21").Value = 399
$ws.Range("L21").Value = 1197
$ws.Range("N21").Value = -1543

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2849707.8
$ws.Range("I113").Value = 6173578.5
$ws.Range("K113").Value = 18520735.5
$ws.Range("M113").Value = -18518565.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 299996
$ws.Range("I128").Value = 299996
$ws.Range("K128").Value = 899988
$ws.Range("M128").Value = -895008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3302.975
$ws.Range("I97").Value = 3988.1292
$ws.Range("K97").Value = 3988.1292
$ws.Range("M97").Value = -3492.1292

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1837301.6
$ws.Range("J132").Value = 337333
$ws.Range("L132").Value = 1011999
$ws.Range("N132").Value = -1017059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2303.75
$ws.Range("I22").Value = 2271.6667
$ws.Range("J22").Value = 2400
$ws.Range("K22").Value = 2271.6667
$ws.Range("L22").Value = 2400
$ws.Range("M22").Value = -1976.6667
$ws.Range("N22").Value = -2990

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 327437.25
$ws.Range("I23").Value = 509124.75
$ws.Range("J23").Value = 145749.75
$ws.Range("K23").Value = 509124.75
$ws.Range("L23").Value = 145749.75
$ws.Range("M23").Value = -508894.75
$ws.Range("N23").Value = -146209.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2303.75
$ws.Range("I27").Value = 2271.6667
$ws.Range("J27").Value = 2400
$ws.Range("K27").Value = 2271.6667
$ws.Range("L27").Value = 2400
$ws.Range("M27").Value = -2164.6667
$ws.Range("N27").Value = -2614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 151258
$ws.Range("I40").Value = 174817.67
$ws.Range("K40").Value = 174817.67
$ws.Range("M40").Value = -174681.67

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4899.9443
$ws.Range("I46").Value = 4319.9
$ws.Range("J46").Value = 5625
$ws.Range("K46").Value = 4319.9
$ws.Range("L46").Value = 5625
$ws.Range("M46").Value = -4131.9
$ws.Range("N46").Value = -6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2527.5715
$ws.Range("I93").Value = 2726
$ws.Range("K93").Value = 2726
$ws.Range("M93").Value = -1478

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 693673.0600000001
$ws.Range("I122").Value = 479814.62
$ws.Range("K122").Value = 1439443.86
$ws.Range("M122").Value = -1436993.86

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 68500
$ws.Range("J140").Value = 68500
$ws.Range("L140").Value = 68500
$ws.Range("N140").Value = -78860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16183.934
$ws.Range("I81").Value = 2090.125
$ws.Range("J81").Value = 32291.143
$ws.Range("K81").Value = 4180.25
$ws.Range("L81").Value = 64582.286
$ws.Range("M81").Value = -3119.25
$ws.Range("N81").Value = -66704.28599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 16183.934
$ws.Range("I84").Value = 2090.125
$ws.Range("J84").Value = 32291.143
$ws.Range("K84").Value = 20901.25
$ws.Range("L84").Value = 322911.43
$ws.Range("M84").Value = -15597.25
$ws.Range("N84").Value = -333519.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 125005544
$ws.Range("J122").Value = 7354.1665
$ws.Range("L122").Value = 22062.4995
$ws.Range("N122").Value = -26962.4995

Write-Output "Applied 32 row updates across 8 sheets (161 cell updates, 1 add, 2 deletes)."
